$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 76

# Column A holds a quarter label formatted like "01-07-2021". Excel's
# auto-detection would otherwise coerce this into a date serial, so force
# it in as literal text (leading apostrophe), then strip the quote-prefix
# formatting it picks up so the cell keeps the sheet's default style.
$cellA = $ws.Range("A$row")
$cellA.Value = "'01-07-2021"
$cellA.Style = "Normal"

$ws.Range("B$row").Value = 233155
$ws.Range("C$row").Value = 80484
$ws.Range("D$row").Value = 44254
$ws.Range("E$row").Value = 3604
$ws.Range("F$row").Value = 4546
$ws.Range("G$row").Value = 28080
$ws.Range("H$row").Value = 152671
$ws.Range("I$row").Value = 22787
$ws.Range("J$row").Value = 6709
$ws.Range("K$row").Value = 123175
$ws.Range("L$row").Value = 25402
$ws.Range("M$row").Value = 2792
$ws.Range("N$row").Value = 0
$ws.Range("O$row").Value = 0
$ws.Range("P$row").Value = 604
$ws.Range("Q$row").Value = 2188
$ws.Range("R$row").Value = 22610
$ws.Range("S$row").Value = 9416
$ws.Range("T$row").Value = 198
$ws.Range("U$row").Value = 12996
$ws.Range("V$row").Value = 207752
$ws.Range("W$row").Value = 77691
$ws.Range("X$row").Value = 44254
$ws.Range("Y$row").Value = 3604
$ws.Range("Z$row").Value = 3942
$ws.Range("AA$row").Value = 25891
$ws.Range("AB$row").Value = 130061
$ws.Range("AC$row").Value = 13371
$ws.Range("AD$row").Value = 6510
$ws.Range("AE$row").Value = 110179
